# RTFL:Concurrency High Level - add new java library reference rows (104-110)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 104 - "More on InputStream vs Reader"
$ws.Range("A104").Value = 'java'
$ws.Range("B104").Value = 'More on InputStream vs Reader'
$ws.Range("C104").Value = 'Reading the Core Java Book, found out that the InputStream is introduced at 1.0, and Reader are introduced later. Basicly InputStream tackle data at a byte level, while Reader tackle data at double-byte/ aka Unicode level. 
The author suggested that we should try to use Reader/Writer at most case, but use the inputStream/ OutputStream when the former is not applicable.'
$ws.Rows.Item(104).RowHeight = 32.25

# Row 105 - "Concurrency: High Level"
$ws.Range("A105").Value = 'RTFM'
$ws.Range("B105").Value = 'Concurrency: High Level '
$ws.Range("C105").Value = 'Thread and Runnable are in java.lang package. But iits lowlevel, to have higher level abstraction, new classes are introduced in _java.util.concurrent package_ after JDK5. 
At the same time, new concurrent package also add new members to the Java Collections Framework.'
$ws.Rows.Item(105).RowHeight = 32.25

# Row 106 - "Concurrency: High Level- Core Interface"
$ws.Range("A106").Value = 'RTFM'
$ws.Range("B106").Value = 'Concurrency: High Level- Core Interface'
$ws.Range("C106").Value = '# Core Interface Introduction
## Executor interface
Provide the core interface with only one method: execute. (e.g. e.execute(r:Runnable);)
## ExecutorService interface
Extends the Executor interface, provide more action on lifecycle management of the service itself and its thread pool memebrs.
## ScheduledExecutorService interface
Extends the ExecutorService and provide scheduling ability
# 2 main ideas in new concurrent pacakge
1. Runnable類從主體變成客體. 以往是new Thread(r).start(); 現在是: executor.execute(r); 好處是應用這個build設計模式使executor的新功能可重用, 強大!
2. ThreadPool concept is introduced to ExecutorService. Basically, a lot of impl of ExecutorService is thread pools. 這是使用Builder pattern 後才使其成為可能.'
$ws.Rows.Item(106).RowHeight = 32.25

# Row 107 - "Concurency: High Level- Executors Factory"
$ws.Range("A107").Value = 'RTFM'
$ws.Range("B107").Value = 'Concurency: High Level- Executors Factory'
$ws.Range("C107").Value = 'Executors(java.util.concurrency.Executors) is a factory that provides different ExecutorService. They are:
* CachedThreadPool
* FixedThreadPool
* ScheduledThreadPool
* SingleThreadExecutor'
$ws.Rows.Item(107).RowHeight = 32.25

# Rows 108-110 - Fork/Join Framework, Collections Framework members, Atomic package
# (headers entered first, long descriptions filled in afterwards - matches original authoring order)
$ws.Range("A108").Value = 'RTFM'
$ws.Range("B108").Value = 'Concurrency: High Level- Fork/Join Framework'
$ws.Range("A109").Value = 'RTFM'
$ws.Range("B109").Value = 'Concurrency: High Level- new members of Java Collection Framework'
$ws.Range("A110").Value = 'RTFM'
$ws.Range("B110").Value = 'Concurrency: High Level- Atomic package'
$ws.Range("C108").Value = 'Fork/Join Framework is introduced to allow the work-stealing algroithm: one free-of-work thread can steal works from another busy thread. 
To enable this framework, refactor the code to follow below pseudo code:
```
if (my portion of the work is small enough)
  do the work directly
else
  split my work into two pieces
  invoke the two pieces and wait for the results
```
For detail usage, reference the manual'
$ws.Range("C109").Value = 'New members are:
* BlockingQueue: defines a first-in-first-out data structure that blocks or times out when you attempt to add to a full queue, or retrieve from an empty queue.
* ConcurrentMap: a subinterface of java.util.Map that defines useful atomic operations. These operations remove or replace a key-value pair only if the key is present, or add a key-value pair only if the key is absent. Making these operations atomic helps avoid synchronization. The standard general-purpose implementation of ConcurrentMap is ConcurrentHashMap, which is a concurrent analog of HashMap.
* ConcurrentNavigableMap: a subinterface of ConcurrentMap that supports approximate matches. The standard general-purpose implementation of ConcurrentNavigableMap is ConcurrentSkipListMap, which is a concurrent analog of TreeMap.'
$ws.Range("C110").Value = 'Atomic package (java.util.concurrency.atomic) provides handy utils for atomic operation. In previous session, the integer counter method increment() adds the **synchronize** keyword to atomicize itself. 
```
class SynchronizedCounter {
    private int c = 0;
    public synchronized void increment() {   c++;   }
    public synchronized void decrement() {   c--;    }
    public synchronized int value() {    return c;    }
}
```
With atomic package, this turns into:
```
import java.util.concurrent.atomic.AtomicInteger;
class AtomicCounter {
    private AtomicInteger c = new AtomicInteger(0);
    public void increment() {        c.incrementAndGet();    }
    public void decrement() {        c.decrementAndGet();    }
    public int value() {        return c.get();    }
}
```'
$ws.Rows.Item(108).RowHeight = 32.25
$ws.Rows.Item(109).RowHeight = 32.25
$ws.Rows.Item(110).RowHeight = 32.25

# Update the view: scroll so row 103 is at top, select C111 (first empty cell below the new data)
$excel.ActiveWindow.ScrollRow = 103
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C111").Select()

